# chore: update Sheets via scheduled runner
# Refreshes market-price derived columns (H:N) across several leve rows
# on the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets. A handful of rows lose
# their NQ-vs-HQ price distinction (the now-stale column is cleared so the
# cell disappears entirely) while the adjacent profit column picks up a
# freshly computed value.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 163.66667
$ws.Range("I18").Value = 163.66667
$ws.Range("K18").Value = 163.66667
$ws.Range("M18").Value = 120.33333
$ws.Range("H33").Value = 4349.3335
$ws.Range("I33").Value = 6799.4
$ws.Range("J33").Value = 265.8889
$ws.Range("K33").Value = 6799.4
$ws.Range("L33").Value = 265.8889
$ws.Range("M33").Value = -6570.4
$ws.Range("N33").Value = -723.8888999999999
$ws.Range("H40").Value = 2027.6471
$ws.Range("I40").Value = 1887.2
$ws.Range("J40").Value = 2228.2856
$ws.Range("K40").Value = 1887.2
$ws.Range("L40").Value = 2228.2856
$ws.Range("M40").Value = -1712.2
$ws.Range("N40").Value = -2578.2856
$ws.Range("H51").Value = 5894.75
$ws.Range("J51").Value = 5894.75
$ws.Range("L51").Value = 5894.75
$ws.Range("N51").Value = -6862.75
$ws.Range("H53").Value = 1000
$ws.Range("I53").Value = 1000
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 1000
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = ""
$ws.Range("N53").Value = -363
$ws.Range("H87").Value = 19894.195
$ws.Range("J87").Value = 19894.195
$ws.Range("L87").Value = 19894.195
$ws.Range("N87").Value = -22390.195
$ws.Range("H90").Value = 19894.195
$ws.Range("J90").Value = 19894.195
$ws.Range("L90").Value = 59682.585
$ws.Range("N90").Value = -72162.58499999999
$ws.Range("H101").Value = 30303456
$ws.Range("I101").Value = 33333682
$ws.Range("J101").Value = 1185
$ws.Range("K101").Value = 100001046
$ws.Range("L101").Value = 3555
$ws.Range("M101").Value = -99999424
$ws.Range("N101").Value = -6799
$ws.Range("H129").Value = 807.4194
$ws.Range("J129").Value = 902.0833
$ws.Range("L129").Value = 2706.2499
$ws.Range("N129").Value = -12706.2499
$ws.Range("H132").Value = 1129.4193
$ws.Range("I132").Value = 761.7692
$ws.Range("K132").Value = 2285.3076
$ws.Range("M132").Value = 244.6923999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6878.761
$ws.Range("I32").Value = 5346.8545
$ws.Range("K32").Value = 5346.8545
$ws.Range("M32").Value = -5059.8545
$ws.Range("H97").Value = 900
$ws.Range("I97").Value = 872.7273
$ws.Range("K97").Value = 872.7273
$ws.Range("M97").Value = -376.7273
$ws.Range("H102").Value = 2107.682
$ws.Range("I102").Value = 2088.45
$ws.Range("K102").Value = 2088.45
$ws.Range("M102").Value = -466.4499999999998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 213.33333
$ws.Range("I22").Value = 213.33333
$ws.Range("K22").Value = 213.33333
$ws.Range("M22").Value = -40.33332999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2583.2856
$ws.Range("I105").Value = 3220.75
$ws.Range("K105").Value = 3220.75
$ws.Range("M105").Value = -1473.75
$ws.Range("H141").Value = 39199.832
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 39199.832
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = ""
$ws.Range("M141").Value = 39199.832
$ws.Range("N141").Value = -49559.832

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1040055.94
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 1126643.9
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 3379931.7
$ws.Range("M107").Value = -1080
$ws.Range("N107").Value = -3383771.7
$ws.Range("H131").Value = 817.79
$ws.Range("I131").Value = 533.8
$ws.Range("J131").Value = 832.7368
$ws.Range("K131").Value = 1601.4
$ws.Range("L131").Value = 2498.2104
$ws.Range("M131").Value = 3438.6
$ws.Range("N131").Value = -12578.2104
$ws.Range("H137").Value = 2432.7715
$ws.Range("I137").Value = 1973.1111
$ws.Range("J137").Value = 2591.8845
$ws.Range("K137").Value = 5919.3333
$ws.Range("L137").Value = 7775.6535
$ws.Range("M137").Value = -819.3333000000002
$ws.Range("N137").Value = -17975.6535

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7719.727
$ws.Range("I70").Value = 5002.6665
$ws.Range("J70").Value = 8738.625
$ws.Range("K70").Value = 5002.6665
$ws.Range("L70").Value = 8738.625
$ws.Range("M70").Value = -4732.6665
$ws.Range("N70").Value = -9278.625
$ws.Range("H73").Value = 7719.727
$ws.Range("I73").Value = 5002.6665
$ws.Range("J73").Value = 8738.625
$ws.Range("K73").Value = 5002.6665
$ws.Range("L73").Value = 8738.625
$ws.Range("M73").Value = -4066.6665
$ws.Range("N73").Value = -10610.625
$ws.Range("H122").Value = 1841.3928
$ws.Range("I122").Value = 1415.1052
$ws.Range("J122").Value = 2741.3333
$ws.Range("K122").Value = 4245.3156
$ws.Range("L122").Value = 8223.999899999999
$ws.Range("M122").Value = -1795.3156
$ws.Range("N122").Value = -13123.9999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 798
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 798
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = ""
$ws.Range("M22").Value = 798
$ws.Range("N22").Value = -1388
$ws.Range("H27").Value = 798
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 798
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = ""
$ws.Range("M27").Value = 798
$ws.Range("N27").Value = -1012
$ws.Range("H46").Value = 505000.5
$ws.Range("I46").Value = 505000.5
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 505000.5
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = ""
$ws.Range("N46").Value = -504812.5
$ws.Range("H55").Value = 571.4286
$ws.Range("I55").Value = 366.33334
$ws.Range("J55").Value = 725.25
$ws.Range("K55").Value = 366.33334
$ws.Range("L55").Value = 725.25
$ws.Range("M55").Value = -193.33334
$ws.Range("N55").Value = -1071.25
$ws.Range("H68").Value = 1732.6154
$ws.Range("I68").Value = 1099
$ws.Range("J68").Value = 1847.8182
$ws.Range("K68").Value = 1099
$ws.Range("L68").Value = 1847.8182
$ws.Range("M68").Value = -350
$ws.Range("N68").Value = -3345.8182
$ws.Range("H71").Value = 1732.6154
$ws.Range("I71").Value = 1099
$ws.Range("J71").Value = 1847.8182
$ws.Range("K71").Value = 5495
$ws.Range("L71").Value = 9239.091
$ws.Range("M71").Value = -1751
$ws.Range("N71").Value = -16727.091
$ws.Range("H133").Value = 36537.8
$ws.Range("J133").Value = 36537.8
$ws.Range("L133").Value = 36537.8
$ws.Range("N133").Value = -41597.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 40495.2
$ws.Range("J46").Value = 40495.2
$ws.Range("L46").Value = 40495.2
$ws.Range("N46").Value = -40957.2
$ws.Range("H81").Value = 2060.6316
$ws.Range("I81").Value = 1944
$ws.Range("J81").Value = 2145.4546
$ws.Range("K81").Value = 3888
$ws.Range("L81").Value = 4290.9092
$ws.Range("M81").Value = -2827
$ws.Range("N81").Value = -6412.9092
$ws.Range("H84").Value = 2060.6316
$ws.Range("I84").Value = 1944
$ws.Range("J84").Value = 2145.4546
$ws.Range("K84").Value = 19440
$ws.Range("L84").Value = 21454.546
$ws.Range("M84").Value = -14136
$ws.Range("N84").Value = -32062.546
$ws.Range("H134").Value = 40495.2
$ws.Range("J134").Value = 40495.2
$ws.Range("L134").Value = 121485.6
$ws.Range("N134").Value = -126555.6
